$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, copy formatting for column A from row 85 (existing formatted row) down to new rows 86-95
$ws.Range("A85").Copy()
$ws.Range("A86:A95").PasteSpecial(-4122)

# Row 80
$ws.Cells.Item(80, 1).Value = 78
$ws.Cells.Item(80, 2).Value = 6.52303547698906
$ws.Cells.Item(80, 3).Value = 11729.4
$ws.Cells.Item(80, 4).Value = 0.06017257123245481
$ws.Cells.Item(80, 5).Value = 106.8
$ws.Cells.Item(80, 6).Value = 284
$ws.Cells.Item(80, 7).Value = "BigFather Rengar"
$ws.Cells.Item(80, 8).Value = "SOLO"
$ws.Cells.Item(80, 9).Value = 0.1657745401248908
$ws.Cells.Item(80, 10).Value = 13
$ws.Cells.Item(80, 11).Value = 0.007222150099747813

# Row 81
$ws.Cells.Item(81, 1).Value = 79
$ws.Cells.Item(81, 2).Value = 3.005318138651472
$ws.Cells.Item(81, 3).Value = 6329.2
$ws.Cells.Item(81, 4).Value = 0.004273504273504274
$ws.Cells.Item(81, 5).Value = 9
$ws.Cells.Item(81, 6).Value = 67
$ws.Cells.Item(81, 7).Value = "Mr Kayn"
$ws.Cells.Item(81, 8).Value = "DUO_SUPPORT"
$ws.Cells.Item(81, 9).Value = 0.03181386514719848
$ws.Cells.Item(81, 10).Value = 5.6
$ws.Cells.Item(81, 11).Value = 0.002659069325735992

# Row 82
$ws.Cells.Item(82, 1).Value = 80
$ws.Cells.Item(82, 2).Value = 6.257362574388895
$ws.Cells.Item(82, 3).Value = 10753.6
$ws.Cells.Item(82, 4).Value = 0.05203516360295726
$ws.Cells.Item(82, 5).Value = 95
$ws.Cells.Item(82, 6).Value = 393.2
$ws.Cells.Item(82, 7).Value = "Booogeyman"
$ws.Cells.Item(82, 8).Value = "DUO_CARRY"
$ws.Cells.Item(82, 9).Value = 0.2543512071886021
$ws.Cells.Item(82, 10).Value = 15.6
$ws.Cells.Item(82, 11).Value = 0.01003492981994257

# Row 83
$ws.Cells.Item(83, 1).Value = 81
$ws.Cells.Item(83, 2).Value = 1.249264705882353
$ws.Cells.Item(83, 3).Value = 2038.8
$ws.Cells.Item(83, 4).Value = 0.01409313725490196
$ws.Cells.Item(83, 5).Value = 23
$ws.Cells.Item(83, 6).Value = 138
$ws.Cells.Item(83, 7).Value = "Poppy Gods"
$ws.Cells.Item(83, 8).Value = "DUO_CARRY"
$ws.Cells.Item(83, 9).Value = 0.08455882352941177
$ws.Cells.Item(83, 10).Value = 3
$ws.Cells.Item(83, 11).Value = 0.001838235294117647

# Row 84
$ws.Cells.Item(84, 1).Value = 82
$ws.Cells.Item(84, 2).Value = 7.410958116892823
$ws.Cells.Item(84, 3).Value = 13175.2
$ws.Cells.Item(84, 4).Value = 0.05164043510720542
$ws.Cells.Item(84, 5).Value = 91.2
$ws.Cells.Item(84, 6).Value = 278.8
$ws.Cells.Item(84, 7).Value = "LS DUFFY"
$ws.Cells.Item(84, 8).Value = "SOLO"
$ws.Cells.Item(84, 9).Value = 0.1567269453783251
$ws.Cells.Item(84, 10).Value = 8.199999999999999
$ws.Cells.Item(84, 11).Value = 0.004627791059971394

# Row 85
$ws.Cells.Item(85, 1).Value = 83
$ws.Cells.Item(85, 2).Value = 7.616710469051133
$ws.Cells.Item(85, 3).Value = 13486.4
$ws.Cells.Item(85, 4).Value = 0.07052820193249709
$ws.Cells.Item(85, 5).Value = 112.4
$ws.Cells.Item(85, 6).Value = 601.6
$ws.Cells.Item(85, 7).Value = "BigFather Rengar"
$ws.Cells.Item(85, 8).Value = "SOLO"
$ws.Cells.Item(85, 9).Value = 0.3650378058503125
$ws.Cells.Item(85, 10).Value = 22.6
$ws.Cells.Item(85, 11).Value = 0.01215785495934224

# Row 86
$ws.Cells.Item(86, 1).Value = 84
$ws.Cells.Item(86, 2).Value = 1.709306742640076
$ws.Cells.Item(86, 3).Value = 3599.8
$ws.Cells.Item(86, 4).Value = 0.002089268755935423
$ws.Cells.Item(86, 5).Value = 4.4
$ws.Cells.Item(86, 6).Value = 140
$ws.Cells.Item(86, 7).Value = "19 fotsiny adc"
$ws.Cells.Item(86, 8).Value = "NONE"
$ws.Cells.Item(86, 9).Value = 0.06647673314339982
$ws.Cells.Item(86, 10).Value = 17
$ws.Cells.Item(86, 11).Value = 0.008072174738841406

# Row 87
$ws.Cells.Item(87, 1).Value = 85
$ws.Cells.Item(87, 2).Value = 6.520541888502177
$ws.Cells.Item(87, 3).Value = 11403.8
$ws.Cells.Item(87, 4).Value = 0.06890925201887924
$ws.Cells.Item(87, 5).Value = 104.4
$ws.Cells.Item(87, 6).Value = 568.8
$ws.Cells.Item(87, 7).Value = "BigFather Rengar"
$ws.Cells.Item(87, 8).Value = "SOLO"
$ws.Cells.Item(87, 9).Value = 0.3700377677125054
$ws.Cells.Item(87, 10).Value = 21.4
$ws.Cells.Item(87, 11).Value = 0.01184026237028092

# Row 88
$ws.Cells.Item(88, 1).Value = 86
$ws.Cells.Item(88, 2).Value = 3.113987360476483
$ws.Cells.Item(88, 3).Value = 5111
$ws.Cells.Item(88, 4).Value = 0.03567156940750339
$ws.Cells.Item(88, 5).Value = 57.8
$ws.Cells.Item(88, 6).Value = 178.4
$ws.Cells.Item(88, 7).Value = "KL  S U S A N O "
$ws.Cells.Item(88, 8).Value = "DUO"
$ws.Cells.Item(88, 9).Value = 0.1073600976001576
$ws.Cells.Item(88, 10).Value = 8
$ws.Cells.Item(88, 11).Value = 0.005132192138943827

# Row 89
$ws.Cells.Item(89, 1).Value = 87
$ws.Cells.Item(89, 2).Value = 6.474714026753915
$ws.Cells.Item(89, 3).Value = 8602
$ws.Cells.Item(89, 4).Value = 0.04634044514716931
$ws.Cells.Item(89, 5).Value = 61.2
$ws.Cells.Item(89, 6).Value = 328.6
$ws.Cells.Item(89, 7).Value = "BigFather Rengar"
$ws.Cells.Item(89, 8).Value = "SOLO"
$ws.Cells.Item(89, 9).Value = 0.2473481266503256
$ws.Cells.Item(89, 10).Value = 9.6
$ws.Cells.Item(89, 11).Value = 0.007220967752454279

# Row 90
$ws.Cells.Item(90, 1).Value = 88
$ws.Cells.Item(90, 2).Value = 2.791646791513082
$ws.Cells.Item(90, 3).Value = 3890.6
$ws.Cells.Item(90, 4).Value = 0.01949175076877424
$ws.Cells.Item(90, 5).Value = 28.2
$ws.Cells.Item(90, 6).Value = 221.4
$ws.Cells.Item(90, 7).Value = "Mr Kayn"
$ws.Cells.Item(90, 8).Value = "DUO_SUPPORT"
$ws.Cells.Item(90, 9).Value = 0.1749684198889241
$ws.Cells.Item(90, 10).Value = 3.2
$ws.Cells.Item(90, 11).Value = 0.002463335941977546

# Row 91
$ws.Cells.Item(91, 1).Value = 93
$ws.Cells.Item(91, 2).Value = 3.07100063734863
$ws.Cells.Item(91, 3).Value = 4818.4
$ws.Cells.Item(91, 4).Value = 0.02179732313575526
$ws.Cells.Item(91, 5).Value = 34.2
$ws.Cells.Item(91, 6).Value = 198
$ws.Cells.Item(91, 7).Value = "Halter Penguen"
$ws.Cells.Item(91, 8).Value = "SOLO"
$ws.Cells.Item(91, 9).Value = 0.1261950286806883
$ws.Cells.Item(91, 10).Value = 4.8
$ws.Cells.Item(91, 11).Value = 0.003059273422562141

# Row 92
$ws.Cells.Item(92, 1).Value = 94
$ws.Cells.Item(92, 2).Value = 3.833485927892299
$ws.Cells.Item(92, 3).Value = 5897.2
$ws.Cells.Item(92, 4).Value = 0.02445485819541133
$ws.Cells.Item(92, 5).Value = 45.2
$ws.Cells.Item(92, 6).Value = 142
$ws.Cells.Item(92, 7).Value = "Mrs Máster Yi"
$ws.Cells.Item(92, 8).Value = "SOLO"
$ws.Cells.Item(92, 9).Value = 0.09318200287130592
$ws.Cells.Item(92, 10).Value = 8.6
$ws.Cells.Item(92, 11).Value = 0.005347130048487147

# Row 93
$ws.Cells.Item(93, 1).Value = 95
$ws.Cells.Item(93, 2).Value = 7.829359333360753
$ws.Cells.Item(93, 3).Value = 16197.8
$ws.Cells.Item(93, 4).Value = 0.08013202994381621
$ws.Cells.Item(93, 5).Value = 165
$ws.Cells.Item(93, 6).Value = 430.2
$ws.Cells.Item(93, 7).Value = "Negabrione"
$ws.Cells.Item(93, 8).Value = "SOLO"
$ws.Cells.Item(93, 9).Value = 0.2063731471137939
$ws.Cells.Item(93, 10).Value = 19.6
$ws.Cells.Item(93, 11).Value = 0.009152606233594018

# Row 94
$ws.Cells.Item(94, 1).Value = 96
$ws.Cells.Item(94, 2).Value = 5.64523200920185
$ws.Cells.Item(94, 3).Value = 8954.200000000001
$ws.Cells.Item(94, 4).Value = 0.02991473419815112
$ws.Cells.Item(94, 5).Value = 49.4
$ws.Cells.Item(94, 6).Value = 610.2
$ws.Cells.Item(94, 7).Value = "HornyCorn"
$ws.Cells.Item(94, 8).Value = "DUO"
$ws.Cells.Item(94, 9).Value = 0.3863354585100985
$ws.Cells.Item(94, 10).Value = 15.2
$ws.Cells.Item(94, 11).Value = 0.009556019006126104

# Row 95
$ws.Cells.Item(95, 1).Value = 97
$ws.Cells.Item(95, 2).Value = 3.652230980643308
$ws.Cells.Item(95, 3).Value = 5297.6
$ws.Cells.Item(95, 4).Value = 0.03020442432290065
$ws.Cells.Item(95, 5).Value = 45.8
$ws.Cells.Item(95, 6).Value = 292
$ws.Cells.Item(95, 7).Value = "MagusApex"
$ws.Cells.Item(95, 8).Value = "DUO"
$ws.Cells.Item(95, 9).Value = 0.1968040937674682
$ws.Cells.Item(95, 10).Value = 5.4
$ws.Cells.Item(95, 11).Value = 0.00344223402841273

